$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: round the Ost/Nord coordinates, drop the empty-looking time cells ---
$ws.Range("Q2").Value = 625863
$ws.Range("R2").Value = 7022772
$ws.Range("Z2").Value = $null
$ws.Range("AB2").Value = $null

# --- Row 3: was the "Norrlandslav" record, becomes the "Bollvitmossa" record ---
$ws.Range("A3").Value = 112094770
$ws.Range("B3").Value = 93881
$ws.Range("E3").Value = 2869
$ws.Range("F3").Value = "Bollvitmossa"
$ws.Range("G3").Value = "Sphagnum wulfianum"
$ws.Range("H3").Value = "Girg."
$ws.Range("Q3").Value = 625863
$ws.Range("R3").Value = 7022772
$ws.Range("Z3").Value = $null
$ws.Range("AB3").Value = $null

# --- Row 4: was the "Bollvitmossa" record, becomes the "Norrlandslav" record ---
$ws.Range("A4").Value = 112094769
$ws.Range("B4").Value = 78604
$ws.Range("E4").Value = 6461
$ws.Range("F4").Value = "Norrlandslav"
$ws.Range("G4").Value = "Nephroma arcticum"
$ws.Range("H4").Value = "(L.) Torss."
$ws.Range("Q4").Value = 625863
$ws.Range("R4").Value = 7022772
$ws.Range("Z4").Value = $null
$ws.Range("AB4").Value = $null
